$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header label
$ws.Range("N1").Value = "Excess +2 to +8"

# Row 2
$ws.Range("B2").Value = 8
$ws.Range("C2").Value = "central store"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "sa"
$ws.Range("F2").Value = "--"
$ws.Range("G2").Value = 200
$ws.Range("H2").Value = 250

# Row 3
$ws.Range("B3").Value = 9
$ws.Range("C3").Value = "test1"
$ws.Range("D3").Value = ""
$ws.Range("F3").Value = "central store"
$ws.Range("G3").Value = 750
$ws.Range("H3").Value = 700

# Row 4
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = "test2"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "sdadads"
$ws.Range("G4").Value = 150
$ws.Range("H4").Value = 100

# Row 5
$ws.Range("B5").Value = 11
$ws.Range("C5").Value = "test3"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "asdad"
$ws.Range("F5").Value = "test2"
$ws.Range("G5").Value = 400
$ws.Range("H5").Value = 1200

# Row 6
$ws.Range("B6").Value = 12
$ws.Range("C6").Value = "mamad"
$ws.Range("D6").Value = ""
$ws.Range("F6").Value = "test3"
$ws.Range("G6").Value = 600
$ws.Range("H6").Value = 600
